$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.472.48'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").Value = '1.830.84'
$ws.Range("E3").Value = '  +2.06%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '316.63'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = '0.5305'
$ws.Range("E7").Value = '  -0.69%  '
$ws.Range("D8").Value = '0.4067'
$ws.Range("E8").Value = '  +7.89%  '
$ws.Range("D9").Value = '0.07571'
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("D10").Value = '41.99'
$ws.Range("E10").Value = '  +0.54%  '
$ws.Range("D11").Value = '1.110'
$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("D12").Value = '6.337'
$ws.Range("E12").Value = '  +3.37%  '
$ws.Range("D13").Value = '1.001'
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("D14").Value = '7.581'
$ws.Range("E14").Value = '  +4.06%  '
$ws.Range("D15").Value = '20.85'
$ws.Range("E15").Value = '  +1.33%  '
$ws.Range("D16").Value = '1.833.67'
$ws.Range("E16").Value = '  +1.92%  '
$ws.Range("D17").Value = '89.56'
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("D18").Value = '0.00001074'
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("D19").Value = '0.06615'
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("D20").Value = '17.55'
$ws.Range("E20").Value = '  +1.21%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").Value = '6.071'
$ws.Range("E22").Value = '  +1.59%  '
$ws.Range("D23").Value = '28.505.88'
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("D24").Value = '11.30'
$ws.Range("E24").Value = '  +2.03%  '
$ws.Range("D25").Value = '2.136'
$ws.Range("E25").Value = '  +2.45%  '
$ws.Range("D26").Value = '2.470'
$ws.Range("E26").Value = '  +8.07%  '
$ws.Range("D27").Value = '157.04'
$ws.Range("E27").Value = '  -1.39%  '
$ws.Range("D28").Value = '20.58'
$ws.Range("E28").Value = '  +0.96%  '
$ws.Range("D29").Value = '2.045.30'
$ws.Range("E29").Value = '  +2.16%  '
$ws.Range("D30").Value = '123.66'
$ws.Range("E30").Value = '  +1.14%  '
$ws.Range("D31").Value = '1.122'
$ws.Range("E31").Value = '  +1.73%  '
$ws.Range("D32").Value = '0.1091'
$ws.Range("E32").Value = '  +3.90%  '
$ws.Range("D33").Value = '5.688'
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("D34").Value = '3.656'
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = '0.07146'
$ws.Range("E35").Value = '  +8.99%  '
$ws.Range("D36").Value = '0.2270'
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("D37").Value = '5.264'
$ws.Range("E37").Value = '  +5.06%  '
$ws.Range("D38").Value = '0.02346'
$ws.Range("E38").Value = '  +2.43%  '
$ws.Range("D39").Value = '8.806'
$ws.Range("E39").Value = '  +3.72%  '
$ws.Range("D40").Value = '11.35'
$ws.Range("E40").Value = '  +2.05%  '
$ws.Range("D41").Value = '0.6272'
$ws.Range("E41").Value = '  +1.64%  '
$ws.Range("D42").Value = '1.192'
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").Value = '1.408'
$ws.Range("E44").Value = '  -3.13%  '
$ws.Range("D45").Value = '13.48'
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("D46").Value = '3.705'
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("D47").Value = '0.5859'
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("D48").Value = '126.10'
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("D49").Value = '1.992'
$ws.Range("E49").Value = '  +2.89%  '
$ws.Range("D50").Value = '1.194'
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("D51").Value = '0.06901'
$ws.Range("E51").Value = '  +0.70%  '
